# The NATMI TPM recompute changed the underlying counts/expression values for
# the Inha -> Tgfbr3 ligand-receptor pair, and added three new sending-cluster
# rows ("MuSCs" as sender) that didn't previously exist in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T, one row per data record (row 1 is the header and is untouched).
$data = @(
    @{ Row=2; Values=@("FAPs","Inha","Tgfbr3","ECs",
            1,0.3333333333333333,0.1623533333333333,0.48706,
            0.2351552968361266,0.2351552968361266,
            3,1,13.10121233333333,39.303637,
            0.1081423012186565,0.1081423012186565,
            2.127025493024445,19.14322943722,
            0.02543023494361497,0.02543023494361498) },
    @{ Row=3; Values=@("FAPs","Inha","Tgfbr3","FAPs",
            1,0.3333333333333333,0.1623533333333333,0.48706,
            0.2351552968361266,0.2351552968361266,
            3,1,69.67747766666666,209.032433,
            0.5751439322003361,0.5751439322003362,
            11.31237075744222,101.81133681698,
            0.1352481421000671,0.1352481421000671) },
    @{ Row=4; Values=@("FAPs","Inha","Tgfbr3","MuSCs",
            1,0.3333333333333333,0.1623533333333333,0.48706,
            0.2351552968361266,0.2351552968361266,
            3,1,38.36920666666666,115.10762,
            0.3167137665810073,0.3167137665810074,
            6.229368599688888,56.0643173972,
            0.07447691979244447,0.07447691979244452) },
    @{ Row=5; Values=@("MuSCs","Inha","Tgfbr3","ECs",
            2,0.6666666666666666,0.5280556666666667,1.584167,
            0.7648447031638734,0.7648447031638734,
            3,1,13.10121233333333,39.303637,
            0.1081423012186565,0.1081423012186565,
            6.91816941281989,62.26352471537901,
            0.08271206627504149,0.0827120662750415) },
    @{ Row=6; Values=@("MuSCs","Inha","Tgfbr3","FAPs",
            2,0.6666666666666666,0.5280556666666667,1.584167,
            0.7648447031638734,0.7648447031638734,
            3,1,69.67747766666666,209.032433,
            0.5751439322003361,0.5751439322003362,
            36.79358692092345,331.142282288311,
            0.439895790100269,0.4398957901002691) },
    @{ Row=7; Values=@("MuSCs","Inha","Tgfbr3","MuSCs",
            2,0.6666666666666666,0.5280556666666667,1.584167,
            0.7648447031638734,0.7648447031638734,
            3,1,38.36920666666666,115.10762,
            0.3167137665810073,0.3167137665810074,
            20.26107700583778,182.34969305254,
            0.2422368467885628,0.2422368467885629) }
)

foreach ($record in $data) {
    $r = $record.Row
    $col = 1
    foreach ($v in $record.Values) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}
